$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: add new row 9 ---
$logs.Cells.Item(9, 1).Value = "Probleem met inloggen"
$logs.Cells.Item(9, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(9, 3).Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$logs.Cells.Item(9, 4).Value = "IT / Technisch probleem"
$logs.Cells.Item(9, 6).Value = "2025-06-20 11:30:25"
$logs.Cells.Item(9, 7).Value = "Nee"

# --- Dashboard sheet: add new row 7 ---
$dash.Cells.Item(7, 1).Value = "IT / Technisch probleem"
$dash.Cells.Item(7, 2).Value = 1

# --- Extend conditional formatting ranges on Logs sheet ---
$fcsD = $logs.Range("D2:D8").FormatConditions
$fcsD.Item(1).ModifyAppliesToRange($logs.Range("D2:D9"))

$fcsG = $logs.Range("G2:G8").FormatConditions
$fcsG.Item(1).ModifyAppliesToRange($logs.Range("G2:G9"))

# --- Update chart series range to include new Dashboard row ---
$co = $dash.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()
$s1 = $sc.Item(1)
$s1.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$7,'Dashboard'!`$B`$2:`$B`$7,1)"
